# 2020_12_12 merge mobile qpyton project directory
#
# Inserts two new command rows ("Unmodifying a Modified File" / git checkout,
# and the git-restore equivalent) right above the existing
# "git ls-files | wc -l" sub-header block (which used to live at row 18 and
# now moves to row 20), re-points the hyperlinks that fall below the
# insertion point, and refreshes the sheet view (top-left cell + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: push rows 18-30 down to 20-32 ------------------------
$ws.Rows.Item(18).Resize(2).Insert()

# New rows should look like the other plain "command / description" rows
# (no custom style / customFormat), so strip whatever got inherited from
# row 17 during the insert.
$ws.Range('A18:C19').ClearFormats()

# --- 2. Fill in the two new rows -----------------------------------------
# Write column B (title) before column A (command) so the new shared
# strings land in the same order as the source workbook.
$ws.Range('B18').Value = 'Unmodifying a Modified File'
$ws.Range('A18').Value = '$ git checkout -- CONTRIBUTING.md'

$ws.Range('B19').Value = 'Unmodifying a Modified File with git restore'
$ws.Range('A19').Value = '$ git restore CONTRIBUTING.md'

# Column C is unused on these rows - make sure no stray empty cell remains.
$ws.Range('C18:C19').ClearContents()

# --- 3. Re-create the hyperlinks at their shifted locations ---------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('C29'), 'https://docs.github.com/en/free-pro-team@latest/github/authenticating-to-github/adding-a-new-ssh-key-to-your-github-account')
$ws.Hyperlinks.Add($ws.Range('C23'), 'https://docs.github.com/en/free-pro-team@latest/github/authenticating-to-github/checking-for-existing-ssh-keys')
$ws.Hyperlinks.Add($ws.Range('C25'), 'https://docs.github.com/en/free-pro-team@latest/github/authenticating-to-github/generating-a-new-ssh-key-and-adding-it-to-the-ssh-agent')
$ws.Hyperlinks.Add($ws.Range('C31'), 'https://docs.github.com/en/free-pro-team@latest/github/authenticating-to-github/testing-your-ssh-connection')
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://git-scm.com/book/en/v2/Git-on-the-Server-Getting-Git-on-a-Server')
$ws.Hyperlinks.Add($ws.Range('C7'), 'https://docs.github.com/en/free-pro-team@latest/github/importing-your-projects-to-github/adding-an-existing-project-to-github-using-the-command-line')
$ws.Hyperlinks.Add($ws.Range('C17'), 'https://backlog.com/git-tutorial/branching/switch-branch/')
$ws.Hyperlinks.Add($ws.Range('C20'), 'https://stackoverflow.com/questions/9468970/how-to-get-a-count-of-all-the-files-in-a-git-repository')

# --- 4. Refresh the view: scroll position + active selection --------------
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range('A27').Select()
